$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update swapped/rotated match rows (columns F:V) ---
# Row 76
$ws.Cells.Item(76, 6).Value = 'Pro Patria'
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 'Mantova'
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 3.38
$ws.Cells.Item(76, 11).Value = '12/10/2023 08:12'
$ws.Cells.Item(76, 12).Value = 4.42
$ws.Cells.Item(76, 13).Value = '14/10/2023 13:51'
$ws.Cells.Item(76, 14).Value = 2.93
$ws.Cells.Item(76, 15).Value = '12/10/2023 08:12'
$ws.Cells.Item(76, 16).Value = 3.21
$ws.Cells.Item(76, 17).Value = '14/10/2023 13:53'
$ws.Cells.Item(76, 18).Value = 2.15
$ws.Cells.Item(76, 19).Value = '12/10/2023 08:12'
$ws.Cells.Item(76, 20).Value = 1.93
$ws.Cells.Item(76, 21).Value = '14/10/2023 13:53'
$ws.Cells.Item(76, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/pro-patria-mantova/UZtBbxWE/'

# Row 77
$ws.Cells.Item(77, 6).Value = 'Pergolettese'
$ws.Cells.Item(77, 7).Value = 2
$ws.Cells.Item(77, 8).Value = 'Trento'
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 2.18
$ws.Cells.Item(77, 11).Value = '12/10/2023 08:12'
$ws.Cells.Item(77, 12).Value = 2.47
$ws.Cells.Item(77, 13).Value = '14/10/2023 13:51'
$ws.Cells.Item(77, 14).Value = 2.9
$ws.Cells.Item(77, 15).Value = '12/10/2023 08:12'
$ws.Cells.Item(77, 16).Value = 2.87
$ws.Cells.Item(77, 17).Value = '14/10/2023 13:51'
$ws.Cells.Item(77, 18).Value = 3.36
$ws.Cells.Item(77, 19).Value = '12/10/2023 08:12'
$ws.Cells.Item(77, 20).Value = 3.29
$ws.Cells.Item(77, 21).Value = '14/10/2023 13:51'
$ws.Cells.Item(77, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/pergolettese-trento/00x7aIG8/'

# Row 101
$ws.Cells.Item(101, 6).Value = 'AlbinoLeffe'
$ws.Cells.Item(101, 7).Value = 1
$ws.Cells.Item(101, 8).Value = 'Arzignano'
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 2.33
$ws.Cells.Item(101, 11).Value = '27/10/2023 02:42'
$ws.Cells.Item(101, 12).Value = 3.05
$ws.Cells.Item(101, 13).Value = '28/10/2023 16:08'
$ws.Cells.Item(101, 14).Value = 2.88
$ws.Cells.Item(101, 15).Value = '27/10/2023 02:42'
$ws.Cells.Item(101, 16).Value = 2.75
$ws.Cells.Item(101, 17).Value = '28/10/2023 16:08'
$ws.Cells.Item(101, 18).Value = 3.07
$ws.Cells.Item(101, 19).Value = '27/10/2023 02:42'
$ws.Cells.Item(101, 20).Value = 2.74
$ws.Cells.Item(101, 21).Value = '28/10/2023 16:08'
$ws.Cells.Item(101, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/albinoleffe-arzignano/l4c544b5/'

# Row 102
$ws.Cells.Item(102, 6).Value = 'Virtus Verona'
$ws.Cells.Item(102, 7).Value = 0
$ws.Cells.Item(102, 8).Value = 'Mantova'
$ws.Cells.Item(102, 9).Value = 2
$ws.Cells.Item(102, 10).Value = 2.21
$ws.Cells.Item(102, 11).Value = '27/10/2023 02:42'
$ws.Cells.Item(102, 12).Value = 2.72
$ws.Cells.Item(102, 13).Value = '28/10/2023 16:07'
$ws.Cells.Item(102, 14).Value = 2.88
$ws.Cells.Item(102, 15).Value = '27/10/2023 02:42'
$ws.Cells.Item(102, 16).Value = 3.09
$ws.Cells.Item(102, 17).Value = '28/10/2023 16:07'
$ws.Cells.Item(102, 18).Value = 3.31
$ws.Cells.Item(102, 19).Value = '27/10/2023 02:42'
$ws.Cells.Item(102, 20).Value = 2.72
$ws.Cells.Item(102, 21).Value = '28/10/2023 16:07'
$ws.Cells.Item(102, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/virtus-verona-mantova/ObXsINS4/'

# Row 104
$ws.Cells.Item(104, 6).Value = 'Giana Erminio'
$ws.Cells.Item(104, 7).Value = 1
$ws.Cells.Item(104, 8).Value = 'Pro Sesto'
$ws.Cells.Item(104, 9).Value = 1
$ws.Cells.Item(104, 10).Value = 2.15
$ws.Cells.Item(104, 11).Value = '27/10/2023 02:42'
$ws.Cells.Item(104, 12).Value = 2.17
$ws.Cells.Item(104, 13).Value = '28/10/2023 18:29'
$ws.Cells.Item(104, 14).Value = 2.93
$ws.Cells.Item(104, 15).Value = '27/10/2023 02:42'
$ws.Cells.Item(104, 16).Value = 3.16
$ws.Cells.Item(104, 17).Value = '28/10/2023 18:29'
$ws.Cells.Item(104, 18).Value = 3.38
$ws.Cells.Item(104, 19).Value = '27/10/2023 02:42'
$ws.Cells.Item(104, 20).Value = 3.57
$ws.Cells.Item(104, 21).Value = '28/10/2023 18:29'
$ws.Cells.Item(104, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/giana-erminio-pro-sesto/ruC83ODB/'

# Row 105
$ws.Cells.Item(105, 6).Value = 'Triestina'
$ws.Cells.Item(105, 7).Value = 2
$ws.Cells.Item(105, 8).Value = 'Fiorenzuola'
$ws.Cells.Item(105, 9).Value = 1
$ws.Cells.Item(105, 10).Value = 1.57
$ws.Cells.Item(105, 11).Value = '27/10/2023 02:42'
$ws.Cells.Item(105, 12).Value = 1.36
$ws.Cells.Item(105, 13).Value = '28/10/2023 18:07'
$ws.Cells.Item(105, 14).Value = 3.51
$ws.Cells.Item(105, 15).Value = '27/10/2023 02:42'
$ws.Cells.Item(105, 16).Value = 4.91
$ws.Cells.Item(105, 17).Value = '28/10/2023 18:07'
$ws.Cells.Item(105, 18).Value = 5.53
$ws.Cells.Item(105, 19).Value = '27/10/2023 02:42'
$ws.Cells.Item(105, 20).Value = 8.710000000000001
$ws.Cells.Item(105, 21).Value = '28/10/2023 18:07'
$ws.Cells.Item(105, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/triestina-fiorenzuola/UyYwJ3Db/'

# Row 122
$ws.Cells.Item(122, 6).Value = 'Fiorenzuola'
$ws.Cells.Item(122, 7).Value = 2
$ws.Cells.Item(122, 8).Value = 'Mantova'
$ws.Cells.Item(122, 9).Value = 3
$ws.Cells.Item(122, 10).Value = 4.33
$ws.Cells.Item(122, 11).Value = '09/11/2023 09:13'
$ws.Cells.Item(122, 12).Value = 4.91
$ws.Cells.Item(122, 13).Value = '11/11/2023 16:11'
$ws.Cells.Item(122, 14).Value = 3.14
$ws.Cells.Item(122, 15).Value = '09/11/2023 09:13'
$ws.Cells.Item(122, 16).Value = 3.66
$ws.Cells.Item(122, 17).Value = '11/11/2023 16:11'
$ws.Cells.Item(122, 18).Value = 1.81
$ws.Cells.Item(122, 19).Value = '09/11/2023 09:13'
$ws.Cells.Item(122, 20).Value = 1.72
$ws.Cells.Item(122, 21).Value = '11/11/2023 16:11'
$ws.Cells.Item(122, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/fiorenzuola-mantova/UHBQSvsN/'

# Row 123
$ws.Cells.Item(123, 6).Value = 'Triestina'
$ws.Cells.Item(123, 7).Value = 0
$ws.Cells.Item(123, 8).Value = 'Pro Sesto'
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 1.33
$ws.Cells.Item(123, 11).Value = '09/11/2023 09:13'
$ws.Cells.Item(123, 12).Value = 1.35
$ws.Cells.Item(123, 13).Value = '11/11/2023 16:12'
$ws.Cells.Item(123, 14).Value = 4.33
$ws.Cells.Item(123, 15).Value = '09/11/2023 09:13'
$ws.Cells.Item(123, 16).Value = 4.56
$ws.Cells.Item(123, 17).Value = '11/11/2023 16:12'
$ws.Cells.Item(123, 18).Value = 9.43
$ws.Cells.Item(123, 19).Value = '09/11/2023 09:13'
$ws.Cells.Item(123, 20).Value = 10.24
$ws.Cells.Item(123, 21).Value = '11/11/2023 16:12'
$ws.Cells.Item(123, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/triestina-pro-sesto/UVjasvdG/'

# Row 125
$ws.Cells.Item(125, 6).Value = 'L.R. Vicenza'
$ws.Cells.Item(125, 7).Value = 3
$ws.Cells.Item(125, 8).Value = 'Pro Patria'
$ws.Cells.Item(125, 9).Value = 1
$ws.Cells.Item(125, 10).Value = 1.75
$ws.Cells.Item(125, 11).Value = '09/11/2023 09:13'
$ws.Cells.Item(125, 12).Value = 1.59
$ws.Cells.Item(125, 13).Value = '11/11/2023 18:21'
$ws.Cells.Item(125, 14).Value = 3.21
$ws.Cells.Item(125, 15).Value = '09/11/2023 09:13'
$ws.Cells.Item(125, 16).Value = 3.64
$ws.Cells.Item(125, 17).Value = '11/11/2023 18:21'
$ws.Cells.Item(125, 18).Value = 4.62
$ws.Cells.Item(125, 19).Value = '09/11/2023 09:13'
$ws.Cells.Item(125, 20).Value = 6.57
$ws.Cells.Item(125, 21).Value = '11/11/2023 18:21'
$ws.Cells.Item(125, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/vicenza-virtus-pro-patria/l2hcNI4p/'

# Row 126
$ws.Cells.Item(126, 6).Value = 'Trento'
$ws.Cells.Item(126, 7).Value = 1
$ws.Cells.Item(126, 8).Value = 'Arzignano'
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 2.16
$ws.Cells.Item(126, 11).Value = '09/11/2023 09:13'
$ws.Cells.Item(126, 12).Value = 2.46
$ws.Cells.Item(126, 13).Value = '11/11/2023 18:21'
$ws.Cells.Item(126, 14).Value = 2.9
$ws.Cells.Item(126, 15).Value = '09/11/2023 09:13'
$ws.Cells.Item(126, 16).Value = 3.06
$ws.Cells.Item(126, 17).Value = '11/11/2023 18:20'
$ws.Cells.Item(126, 18).Value = 3.41
$ws.Cells.Item(126, 19).Value = '09/11/2023 09:13'
$ws.Cells.Item(126, 20).Value = 3.09
$ws.Cells.Item(126, 21).Value = '11/11/2023 18:21'
$ws.Cells.Item(126, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/trento-arzignano/z9ierKtA/'

# Row 132
$ws.Cells.Item(132, 6).Value = 'Virtus Verona'
$ws.Cells.Item(132, 7).Value = 0
$ws.Cells.Item(132, 8).Value = 'Lumezzane'
$ws.Cells.Item(132, 9).Value = 2
$ws.Cells.Item(132, 10).Value = 1.74
$ws.Cells.Item(132, 11).Value = '16/11/2023 09:12'
$ws.Cells.Item(132, 12).Value = 1.83
$ws.Cells.Item(132, 13).Value = '18/11/2023 13:57'
$ws.Cells.Item(132, 14).Value = 3.27
$ws.Cells.Item(132, 15).Value = '16/11/2023 09:12'
$ws.Cells.Item(132, 16).Value = 3.1
$ws.Cells.Item(132, 17).Value = '18/11/2023 13:57'
$ws.Cells.Item(132, 18).Value = 4.56
$ws.Cells.Item(132, 19).Value = '16/11/2023 09:12'
$ws.Cells.Item(132, 20).Value = 5.31
$ws.Cells.Item(132, 21).Value = '18/11/2023 13:57'
$ws.Cells.Item(132, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/virtus-verona-lumezzane/rsP6wVfr/'

# Row 133
$ws.Cells.Item(133, 6).Value = 'Pro Patria'
$ws.Cells.Item(133, 7).Value = 1
$ws.Cells.Item(133, 8).Value = 'Fiorenzuola'
$ws.Cells.Item(133, 9).Value = 1
$ws.Cells.Item(133, 10).Value = 1.98
$ws.Cells.Item(133, 11).Value = '16/11/2023 09:12'
$ws.Cells.Item(133, 12).Value = 2.03
$ws.Cells.Item(133, 13).Value = '18/11/2023 13:51'
$ws.Cells.Item(133, 14).Value = 3.02
$ws.Cells.Item(133, 15).Value = '16/11/2023 09:12'
$ws.Cells.Item(133, 16).Value = 3.23
$ws.Cells.Item(133, 17).Value = '18/11/2023 13:51'
$ws.Cells.Item(133, 18).Value = 3.77
$ws.Cells.Item(133, 19).Value = '16/11/2023 09:12'
$ws.Cells.Item(133, 20).Value = 3.9
$ws.Cells.Item(133, 21).Value = '18/11/2023 13:51'
$ws.Cells.Item(133, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/pro-patria-fiorenzuola/OhCYewB9/'

# Row 134
$ws.Cells.Item(134, 6).Value = 'Alessandria'
$ws.Cells.Item(134, 7).Value = 1
$ws.Cells.Item(134, 8).Value = 'Giana Erminio'
$ws.Cells.Item(134, 9).Value = 2
$ws.Cells.Item(134, 10).Value = 2.46
$ws.Cells.Item(134, 11).Value = '16/11/2023 09:12'
$ws.Cells.Item(134, 12).Value = 2.98
$ws.Cells.Item(134, 13).Value = '18/11/2023 13:56'
$ws.Cells.Item(134, 14).Value = 2.85
$ws.Cells.Item(134, 15).Value = '16/11/2023 09:12'
$ws.Cells.Item(134, 16).Value = 2.78
$ws.Cells.Item(134, 17).Value = '18/11/2023 13:56'
$ws.Cells.Item(134, 18).Value = 2.91
$ws.Cells.Item(134, 19).Value = '16/11/2023 09:12'
$ws.Cells.Item(134, 20).Value = 2.76
$ws.Cells.Item(134, 21).Value = '18/11/2023 13:56'
$ws.Cells.Item(134, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/alessandria-giana-erminio/8Mk3tbBM/'

# Row 135
$ws.Cells.Item(135, 6).Value = 'Arzignano'
$ws.Cells.Item(135, 7).Value = 1
$ws.Cells.Item(135, 8).Value = 'Renate'
$ws.Cells.Item(135, 9).Value = 1
$ws.Cells.Item(135, 10).Value = 2.31
$ws.Cells.Item(135, 11).Value = '16/11/2023 09:12'
$ws.Cells.Item(135, 12).Value = 3.11
$ws.Cells.Item(135, 13).Value = '18/11/2023 13:56'
$ws.Cells.Item(135, 14).Value = 2.86
$ws.Cells.Item(135, 15).Value = '16/11/2023 09:12'
$ws.Cells.Item(135, 16).Value = 2.95
$ws.Cells.Item(135, 17).Value = '18/11/2023 13:56'
$ws.Cells.Item(135, 18).Value = 3.13
$ws.Cells.Item(135, 19).Value = '16/11/2023 09:12'
$ws.Cells.Item(135, 20).Value = 2.52
$ws.Cells.Item(135, 21).Value = '18/11/2023 13:37'
$ws.Cells.Item(135, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/arzignano-renate/Qcp8uIQS/'

# Row 136
$ws.Cells.Item(136, 6).Value = 'Legnago Salus'
$ws.Cells.Item(136, 7).Value = 1
$ws.Cells.Item(136, 8).Value = 'L.R. Vicenza'
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 4.2
$ws.Cells.Item(136, 11).Value = '16/11/2023 09:12'
$ws.Cells.Item(136, 12).Value = 6.68
$ws.Cells.Item(136, 13).Value = '18/11/2023 13:55'
$ws.Cells.Item(136, 14).Value = 3.21
$ws.Cells.Item(136, 15).Value = '16/11/2023 09:12'
$ws.Cells.Item(136, 16).Value = 3.21
$ws.Cells.Item(136, 17).Value = '18/11/2023 13:55'
$ws.Cells.Item(136, 18).Value = 1.81
$ws.Cells.Item(136, 19).Value = '16/11/2023 09:12'
$ws.Cells.Item(136, 20).Value = 1.68
$ws.Cells.Item(136, 21).Value = '18/11/2023 13:55'
$ws.Cells.Item(136, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/legnago-salus-vicenza-virtus/MTVwzdlj/'

# Row 143
$ws.Cells.Item(143, 6).Value = 'Pro Vercelli'
$ws.Cells.Item(143, 7).Value = 2
$ws.Cells.Item(143, 8).Value = 'Alessandria'
$ws.Cells.Item(143, 9).Value = 0
$ws.Cells.Item(143, 10).Value = 1.69
$ws.Cells.Item(143, 11).Value = '23/11/2023 09:13'
$ws.Cells.Item(143, 12).Value = 1.77
$ws.Cells.Item(143, 13).Value = '25/11/2023 18:15'
$ws.Cells.Item(143, 14).Value = 3.26
$ws.Cells.Item(143, 15).Value = '23/11/2023 09:13'
$ws.Cells.Item(143, 16).Value = 3.39
$ws.Cells.Item(143, 17).Value = '25/11/2023 18:15'
$ws.Cells.Item(143, 18).Value = 4.9
$ws.Cells.Item(143, 19).Value = '23/11/2023 09:13'
$ws.Cells.Item(143, 20).Value = 5.05
$ws.Cells.Item(143, 21).Value = '25/11/2023 18:15'
$ws.Cells.Item(143, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/pro-vercelli-alessandria/j3CDdA1K/'

# Row 144
$ws.Cells.Item(144, 6).Value = 'Renate'
$ws.Cells.Item(144, 7).Value = 1
$ws.Cells.Item(144, 8).Value = 'Virtus Verona'
$ws.Cells.Item(144, 9).Value = 1
$ws.Cells.Item(144, 10).Value = 2.46
$ws.Cells.Item(144, 11).Value = '23/11/2023 09:13'
$ws.Cells.Item(144, 12).Value = 2.51
$ws.Cells.Item(144, 13).Value = '25/11/2023 18:24'
$ws.Cells.Item(144, 14).Value = 2.83
$ws.Cells.Item(144, 15).Value = '23/11/2023 09:13'
$ws.Cells.Item(144, 16).Value = 2.88
$ws.Cells.Item(144, 17).Value = '25/11/2023 18:17'
$ws.Cells.Item(144, 18).Value = 2.93
$ws.Cells.Item(144, 19).Value = '23/11/2023 09:13'
$ws.Cells.Item(144, 20).Value = 3.22
$ws.Cells.Item(144, 21).Value = '25/11/2023 18:24'
$ws.Cells.Item(144, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/renate-virtus-verona/YR1IeUGQ/'

# Row 146
$ws.Cells.Item(146, 6).Value = 'L.R. Vicenza'
$ws.Cells.Item(146, 7).Value = 2
$ws.Cells.Item(146, 8).Value = 'Pro Sesto'
$ws.Cells.Item(146, 9).Value = 0
$ws.Cells.Item(146, 10).Value = 1.34
$ws.Cells.Item(146, 11).Value = '23/11/2023 09:13'
$ws.Cells.Item(146, 12).Value = 1.45
$ws.Cells.Item(146, 13).Value = '25/11/2023 20:36'
$ws.Cells.Item(146, 14).Value = 4.31
$ws.Cells.Item(146, 15).Value = '23/11/2023 09:13'
$ws.Cells.Item(146, 16).Value = 4.08
$ws.Cells.Item(146, 17).Value = '25/11/2023 20:43'
$ws.Cells.Item(146, 18).Value = 9.01
$ws.Cells.Item(146, 19).Value = '23/11/2023 09:13'
$ws.Cells.Item(146, 20).Value = 8.27
$ws.Cells.Item(146, 21).Value = '25/11/2023 20:36'
$ws.Cells.Item(146, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/vicenza-virtus-pro-sesto/UcRIzTv1/'

# Row 147
$ws.Cells.Item(147, 6).Value = 'Legnago Salus'
$ws.Cells.Item(147, 7).Value = 1
$ws.Cells.Item(147, 8).Value = 'Pro Patria'
$ws.Cells.Item(147, 9).Value = 1
$ws.Cells.Item(147, 10).Value = 2.3
$ws.Cells.Item(147, 11).Value = '23/11/2023 21:42'
$ws.Cells.Item(147, 12).Value = 2.39
$ws.Cells.Item(147, 13).Value = '25/11/2023 20:41'
$ws.Cells.Item(147, 14).Value = 2.91
$ws.Cells.Item(147, 15).Value = '23/11/2023 21:42'
$ws.Cells.Item(147, 16).Value = 2.92
$ws.Cells.Item(147, 17).Value = '25/11/2023 20:41'
$ws.Cells.Item(147, 18).Value = 3.2
$ws.Cells.Item(147, 19).Value = '23/11/2023 21:42'
$ws.Cells.Item(147, 20).Value = 3.36
$ws.Cells.Item(147, 21).Value = '25/11/2023 20:41'
$ws.Cells.Item(147, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/legnago-salus-pro-patria/nZLNZnf7/'

# Row 148
$ws.Cells.Item(148, 6).Value = 'Fiorenzuola'
$ws.Cells.Item(148, 7).Value = 1
$ws.Cells.Item(148, 8).Value = 'Giana Erminio'
$ws.Cells.Item(148, 9).Value = 2
$ws.Cells.Item(148, 10).Value = 2.53
$ws.Cells.Item(148, 11).Value = '23/11/2023 09:13'
$ws.Cells.Item(148, 12).Value = 2.75
$ws.Cells.Item(148, 13).Value = '26/11/2023 13:58'
$ws.Cells.Item(148, 14).Value = 2.89
$ws.Cells.Item(148, 15).Value = '23/11/2023 09:13'
$ws.Cells.Item(148, 16).Value = 2.99
$ws.Cells.Item(148, 17).Value = '26/11/2023 13:50'
$ws.Cells.Item(148, 18).Value = 2.77
$ws.Cells.Item(148, 19).Value = '23/11/2023 09:13'
$ws.Cells.Item(148, 20).Value = 2.78
$ws.Cells.Item(148, 21).Value = '26/11/2023 13:58'
$ws.Cells.Item(148, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/fiorenzuola-giana-erminio/j7SEy9Pf/'

# Row 149
$ws.Cells.Item(149, 6).Value = 'Lumezzane'
$ws.Cells.Item(149, 7).Value = 1
$ws.Cells.Item(149, 8).Value = 'Novara'
$ws.Cells.Item(149, 9).Value = 1
$ws.Cells.Item(149, 10).Value = 2.05
$ws.Cells.Item(149, 11).Value = '23/11/2023 09:13'
$ws.Cells.Item(149, 12).Value = 2.33
$ws.Cells.Item(149, 13).Value = '26/11/2023 13:50'
$ws.Cells.Item(149, 14).Value = 2.99
$ws.Cells.Item(149, 15).Value = '23/11/2023 09:13'
$ws.Cells.Item(149, 16).Value = 3.02
$ws.Cells.Item(149, 17).Value = '26/11/2023 13:50'
$ws.Cells.Item(149, 18).Value = 3.59
$ws.Cells.Item(149, 19).Value = '23/11/2023 09:13'
$ws.Cells.Item(149, 20).Value = 3.37
$ws.Cells.Item(149, 21).Value = '26/11/2023 13:50'
$ws.Cells.Item(149, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/lumezzane-novara/QPKRY69D/'

# --- Add new row 150 (copy formatting from row 149 first) ---
$ws.Range("A149:V149").Copy()
$ws.Range("A150:V150").PasteSpecial(-4122)
$ws.Cells.Item(150, 1).Value = 149
$ws.Cells.Item(150, 2).Value = 'italy'
$ws.Cells.Item(150, 3).Value = 'serie-c-group-a'
$ws.Cells.Item(150, 4).Value = '2023-2024'
$ws.Cells.Item(150, 5).Value = 45257.86458333334
$ws.Cells.Item(150, 6).Value = 'Pergolettese'
$ws.Cells.Item(150, 7).Value = 1
$ws.Cells.Item(150, 8).Value = 'Padova'
$ws.Cells.Item(150, 9).Value = 1
$ws.Cells.Item(150, 10).Value = 3.63
$ws.Cells.Item(150, 11).Value = '23/11/2023 09:13'
$ws.Cells.Item(150, 12).Value = 4.5
$ws.Cells.Item(150, 13).Value = '27/11/2023 20:41'
$ws.Cells.Item(150, 14).Value = 3.09
$ws.Cells.Item(150, 15).Value = '23/11/2023 09:13'
$ws.Cells.Item(150, 16).Value = 3.36
$ws.Cells.Item(150, 17).Value = '27/11/2023 20:39'
$ws.Cells.Item(150, 18).Value = 1.99
$ws.Cells.Item(150, 19).Value = '23/11/2023 09:13'
$ws.Cells.Item(150, 20).Value = 1.86
$ws.Cells.Item(150, 21).Value = '27/11/2023 20:41'
$ws.Cells.Item(150, 22).Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/pergolettese-padova/vT89cjnE/'
